# "change tracker to db" — append rows 3-16 pulled from the DB sync, matching
# the existing header/row-2 formatting, and fix up row 2's PRODUCT_ID /
# LINE_CREATION_DATE / Config_error text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 touch-ups -------------------------------------------------
$ws.Cells.Item(2, 6).Value  = "C9410R-96U-BNDL-A"   # F2 PRODUCT_ID (unchanged text, re-set for parity)
$ws.Cells.Item(2, 10).Value = "SUP qty over"         # J2 Config_error (unchanged text, re-set for parity)

# H2 (LINE_CREATION_DATE) must stay the literal text "02/11/2021" rather than
# being auto-parsed into a date serial. Build it via a formula then convert
# that formula result to a static value with PasteSpecial so Excel stores it
# as a shared string instead of re-interpreting the text as a date.
$h2 = $ws.Cells.Item(2, 8)
$h2.Formula = "=""02/11/2021"""
$h2.Copy()
$h2.PasteSpecial(-4163)   # xlPasteValues

# --- New rows 3-16 (DB export rows) -----------------------------------
$newRows = @(
    @{ Row=3;  D="111142195-37"; F="C9410R";            K=44348.42330030093 },
    @{ Row=4;  D="111142195-41"; F="C9410R";            K=44348.42330030093 },
    @{ Row=5;  D="111142195-47"; F="C9410R";            K=44348.42330030093 },
    @{ Row=6;  D="111142918-4";  F="C9407R-96U-BNDL-A"; K=44348.42330030093 },
    @{ Row=7;  D="111142195-54"; F="C9410R";            K=44348.42330030093 },
    @{ Row=8;  D="111142195-59"; F="C9410R";            K=44348.42330030093 },
    @{ Row=9;  D="111142195-73"; F="C9410R";            K=44348.42330030093 },
    @{ Row=10; D="111142195-77"; F="C9410R";            K=44348.42330030093 },
    @{ Row=11; D="111142195-81"; F="C9410R";            K=44348.42330030093 },
    @{ Row=12; D="111142195-85"; F="C9410R";            K=44348.42330030093 },
    @{ Row=13; D="111142195-89"; F="C9410R";            K=44348.42330030093 },
    @{ Row=14; D="111424230-2";  F="C9410R-96U-BNDL-A"; K=44348.42330030093 },
    @{ Row=15; D="111401757-3";  F="C9410R";            K=44348.42330030093 },
    @{ Row=16; D="112466942-1";  F="ISR4331-DNA";       K=44348.42637737964 }
)

# Fill column D (PO_NUMBER) for every new row first, then column F
# (PRODUCT_ID), then column J (Config_error) — matches the order the rows
# were originally populated in.
foreach ($r in $newRows) {
    $ws.Cells.Item($r.Row, 4).Value = $r.D
}
foreach ($r in $newRows) {
    $ws.Cells.Item($r.Row, 6).Value = $r.F
}
foreach ($r in $newRows) {
    $ws.Cells.Item($r.Row, 10).Value = "(user report) test"
}

foreach ($r in $newRows) {
    $row = $r.Row

    $ws.Cells.Item($row, 1).Value = "JMX"   # ORGANIZATION_CODE
    $ws.Cells.Item($row, 2).Value = "UABU"  # BUSINESS_UNIT
    $ws.Cells.Item($row, 5).Value = 0       # OPTION_NUMBER
    $ws.Cells.Item($row, 7).Value = 1       # ORDERED_QUANTITY

    $k = $ws.Cells.Item($row, 11)
    $k.Value = $r.K
    $k.NumberFormat = "YYYY-MM-DD HH:MM:SS"

    # Match row 2's look: bold/bordered/centered ORGANIZATION_CODE cell,
    # everything else plain.
    $ws.Range("A2").Copy()
    $ws.Cells.Item($row, 1).PasteSpecial(-4122)   # xlPasteFormats
}

$null = $ws.Range("A1").Select()
